# B6-PowerPoint.pptx — replay of the "Tue, May 05, 2020  9:09:21 PM" commit.
#
# The commit does two things to the deck:
#   1. Re-styles the three tables (slides 14, 15, 16) from table style
#      {502CBCC7-B206-4623-8CA5-1D2FDDFE635F} to
#      {82DD435E-599B-44C1-B8CE-AC7D54BADD7D}.
#   2. Swaps the presentation's theme colours: the deck's live theme
#      ("Integral" / "Red Violet") is replaced with the plain "Office"
#      colour palette (the colours that used to live, unused, in the
#      deck's secondary theme part).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style swap — one table per slide, always the first shape.
# ---------------------------------------------------------------------
$newTableStyle = "{82DD435E-599B-44C1-B8CE-AC7D54BADD7D}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyle)
    }
}

# ---------------------------------------------------------------------
# 2) Theme colour swap — push the "Office" colour scheme onto the
#    presentation's theme colours (dk1/lt1/dk2/lt2/accent1-6/hlink/
#    folHlink, in that order).
# ---------------------------------------------------------------------
$officeColors = 0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477

$theme = $p.Designs.Item(1).SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}

Write-Host "Applied table style $newTableStyle to slides 14/15/16 and reset theme colors to Office palette."
